$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.191131666666666
$ws.Range("H2").Value = 3.573395
$ws.Range("I2").Value = 0.02720036629735778
$ws.Range("J2").Value = 0.02720036629735778
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 103.4766596666667
$ws.Range("N2").Value = 310.429979
$ws.Range("O2").Value = 0.877785331764719
$ws.Range("P2").Value = 0.8777853317647188
$ws.Range("Q2").Value = 123.2543260898561
$ws.Range("R2").Value = 1109.288934808705
$ws.Range("S2").Value = 0.02387608255444808
$ws.Range("T2").Value = 0.02387608255444807

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.191131666666666
$ws.Range("H3").Value = 3.573395
$ws.Range("I3").Value = 0.02720036629735778
$ws.Range("J3").Value = 0.02720036629735778
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8265796666666668
$ws.Range("N3").Value = 2.479739
$ws.Range("O3").Value = 0.007011818020336602
$ws.Range("P3").Value = 0.0070118180203366
$ws.Range("Q3").Value = 0.9845652159894445
$ws.Range("R3").Value = 8.861086943905
$ws.Range("S3").Value = 0.0001907240185635696
$ws.Range("T3").Value = 0.0001907240185635696

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.191131666666666
$ws.Range("H4").Value = 3.573395
$ws.Range("I4").Value = 0.02720036629735778
$ws.Range("J4").Value = 0.02720036629735778
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.58054833333333
$ws.Range("N4").Value = 40.741645
$ws.Range("O4").Value = 0.1152028502149446
$ws.Range("P4").Value = 0.1152028502149446
$ws.Range("Q4").Value = 16.17622117053055
$ws.Range("R4").Value = 145.585990534775
$ws.Range("S4").Value = 0.003133559724346136
$ws.Range("T4").Value = 0.003133559724346135

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 34.415161
$ws.Range("H5").Value = 103.245483
$ws.Range("I5").Value = 0.7858954736735307
$ws.Range("J5").Value = 0.7858954736735306
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 103.4766596666667
$ws.Range("N5").Value = 310.429979
$ws.Range("O5").Value = 0.877785331764719
$ws.Range("P5").Value = 0.8777853317647188
$ws.Range("Q5").Value = 3561.16590217054
$ws.Range("R5").Value = 32050.49311953486
$ws.Range("S5").Value = 0.6898475190909111
$ws.Range("T5").Value = 0.6898475190909109

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 34.415161
$ws.Range("H6").Value = 103.245483
$ws.Range("I6").Value = 0.7858954736735307
$ws.Range("J6").Value = 0.7858954736735306
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8265796666666668
$ws.Range("N6").Value = 2.479739
$ws.Range("O6").Value = 0.007011818020336602
$ws.Range("P6").Value = 0.0070118180203366
$ws.Range("Q6").Value = 28.44687230765967
$ws.Range("R6").Value = 256.021850768937
$ws.Range("S6").Value = 0.005510556044405032
$ws.Range("T6").Value = 0.00551055604440503

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 34.415161
$ws.Range("H7").Value = 103.245483
$ws.Range("I7").Value = 0.7858954736735307
$ws.Range("J7").Value = 0.7858954736735306
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.58054833333333
$ws.Range("N7").Value = 40.741645
$ws.Range("O7").Value = 0.1152028502149446
$ws.Range("P7").Value = 0.1152028502149446
$ws.Range("Q7").Value = 467.3767573599484
$ws.Range("R7").Value = 4206.390816239536
$ws.Range("S7").Value = 0.09053739853821471
$ws.Range("T7").Value = 0.09053739853821467

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.184723
$ws.Range("H8").Value = 24.554169
$ws.Range("I8").Value = 0.1869041600291116
$ws.Range("J8").Value = 0.1869041600291116
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.4766596666667
$ws.Range("N8").Value = 310.429979
$ws.Range("O8").Value = 0.877785331764719
$ws.Range("P8").Value = 0.8777853317647188
$ws.Range("Q8").Value = 846.927796336939
$ws.Range("R8").Value = 7622.350167032451
$ws.Range("S8").Value = 0.1640617301193599
$ws.Range("T8").Value = 0.1640617301193598

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.184723
$ws.Range("H9").Value = 24.554169
$ws.Range("I9").Value = 0.1869041600291116
$ws.Range("J9").Value = 0.1869041600291116
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8265796666666668
$ws.Range("N9").Value = 2.479739
$ws.Range("O9").Value = 0.007011818020336602
$ws.Range("P9").Value = 0.0070118180203366
$ws.Range("Q9").Value = 6.765325609099001
$ws.Range("R9").Value = 60.88793048189101
$ws.Range("S9").Value = 0.001310537957368001
$ws.Range("T9").Value = 0.001310537957368

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.184723
$ws.Range("H10").Value = 24.554169
$ws.Range("I10").Value = 0.1869041600291116
$ws.Range("J10").Value = 0.1869041600291116
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.58054833333333
$ws.Range("N10").Value = 40.741645
$ws.Range("O10").Value = 0.1152028502149446
$ws.Range("P10").Value = 0.1152028502149446
$ws.Range("Q10").Value = 111.153026296445
$ws.Range("R10").Value = 1000.377236668005
$ws.Range("S10").Value = 0.02153189195238378
$ws.Range("T10").Value = 0.02153189195238377

